$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)
$hdr = $sec.Headers.Item(1)
$rng = $hdr.Range
$ok = $rng.Find.Execute("Oscilloscope")
Write-Host "ok:" $ok "Start:" $rng.Start "End:" $rng.End
$rng.MoveEnd(1, -1)
Write-Host "Start:" $rng.Start "End:" $rng.End "text:[$($rng.Text)]"
$rng.Text = ""
$final = $hdr.Range
Write-Host "Final:" $final.Text
